# correction scrappe Ag LBMA
# Append a fresh batch of "25/05/2023" price rows to every per-metal sheet
# (mirroring the existing row 3/4 values down through row 9, or row 5 for
# the two sheets that only had a single 25/05/2023 sample), and refresh the
# Ag (1AG2) value on the RPA roll-up sheet to match.

$wb = $excel.ActiveWorkbook

function Add-Rows {
    # Positional params only - named binding ("-Foo bar") is unreliable here.
    param($SheetName, $FirstRow, $LastRow, $BValue, $BIsText, $CValue, $DValue)

    $ws = $wb.Worksheets.Item($SheetName)

    for ($r = $FirstRow; $r -le $LastRow; $r++) {
        $ws.Range("A$r").Value = "25/05/2023"

        if ($BIsText) {
            # Force text storage so values like "687,70" or "23,415" are not
            # re-interpreted as numbers (comma grouping) by the assignment.
            $ws.Range("B$r").NumberFormat = "@"
        }
        $ws.Range("B$r").Value = $BValue

        $ws.Range("C$r").Value = $CValue
        $ws.Range("D$r").Value = $DValue
    }
}

# 1AG1 - Ag, KG, EUR : rows 5-9, numeric 1234
Add-Rows "1AG1" 5 9 1234 $false "€" "KG"

# 1AG3 - Ag, KG, EUR : rows 5-9, text "687,70"
Add-Rows "1AG3" 5 9 "687,70" $true "€" "KG"

# 1AU3 - Au, KG, EUR : rows 5-9, numeric 123
Add-Rows "1AU3" 5 9 123 $false "€" "KG"

# 2M37 - 100KG, EUR : rows 5-9, text "718,00"
Add-Rows "2M37" 5 9 "718,00" $true "€" "100KG"

# 3AL1 - Aluminium, TO, USD : rows 5-9, text "2234,00"
Add-Rows "3AL1" 5 9 "2234,00" $true "$" "TO"

# 3CU1 - Copper, TO, USD : rows 5-9, text "7910,00"
Add-Rows "3CU1" 5 9 "7910,00" $true "$" "TO"

# 3CU3 - Copper, 100KG, EUR : rows 5-9, text "873,13"
Add-Rows "3CU3" 5 9 "873,13" $true "€" "100KG"

# 2CUB - Copper, KG, EUR : rows 5-9, text "9,27"
Add-Rows "2CUB" 5 9 "9,27" $true "€" "KG"

# 3NI1 - Nickel, TO, USD : rows 5-9, text "20900,00"
Add-Rows "3NI1" 5 9 "20900,00" $true "$" "TO"

# 3SN1 - Tin, TO, USD : rows 5-9, text "24550,00"
Add-Rows "3SN1" 5 9 "24550,00" $true "$" "TO"

# 1AG2 - Ag, OZ, USD : row 5 only, text "23,415"
Add-Rows "1AG2" 5 5 "23,415" $true "$" "OZ"

# 1AU2 - Au, OZ, USD : row 5 only, text "1969,65"
Add-Rows "1AU2" 5 5 "1969,65" $true "$" "OZ"

# RPA roll-up sheet: refresh the 1AG2 (Ag/OZ) quote to the corrected value
$rpa = $wb.Worksheets.Item("RPA")
$rpa.Range("B12").NumberFormat = "@"
$rpa.Range("B12").Value = "23,415"
